# Generate Report for Handoff
# Appends two new source-file records (1d9b0e82-... and 3c3f2e5e-...) to the
# Overview sheet and to each per-language detail sheet (zh-cn, de-de),
# mirroring the existing "Ready for handoff" row pattern (e.g. the
# 01c3da70-f464-4d60-973a-d00a275bd8ed row) since neither new file has been
# handed back yet.

$wb = $excel.ActiveWorkbook

$records = @(
    @{
        Uuid         = "1d9b0e82-f493-4953-be8e-bf9b78a9a8df"
        MdUrl        = "https://github.com/OpenLocalizationTest/oltest/blob/d9de1136c2a2a229cd87b2f958516c4145a80933/e2e/1d9b0e82-f493-4953-be8e-bf9b78a9a8df.md"
        HandoffDate  = "2016-31-12 18:31:07"
        ZhXlfFile    = "1d9b0e82-f493-4953-be8e-bf9b78a9a8df.ffc11811f9871f582b0d4cf682cbc358c94f41dd.zh-cn.xlf"
        ZhXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/53d75881395e0b17c48027a30078233ddac0ec07/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1d9b0e82-f493-4953-be8e-bf9b78a9a8df.ffc11811f9871f582b0d4cf682cbc358c94f41dd.zh-cn.xlf"
        ZhHandoffDt  = "2016-03-12 18:31:04"
        DeXlfFile    = "1d9b0e82-f493-4953-be8e-bf9b78a9a8df.ffc11811f9871f582b0d4cf682cbc358c94f41dd.de-de.xlf"
        DeXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/17f07534802ef9d1eb6e0f01ccd138bb77467e79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1d9b0e82-f493-4953-be8e-bf9b78a9a8df.ffc11811f9871f582b0d4cf682cbc358c94f41dd.de-de.xlf"
        DeHandoffDt  = "2016-03-12 18:31:07"
    },
    @{
        Uuid         = "3c3f2e5e-9ad4-413f-8618-f6ee6d42b7ff"
        MdUrl        = "https://github.com/OpenLocalizationTest/oltest/blob/28aae8cfdd3d619e990c3a0f53c0a3cb3cea15c4/e2e/3c3f2e5e-9ad4-413f-8618-f6ee6d42b7ff.md"
        HandoffDate  = "2016-31-12 18:31:07"
        ZhXlfFile    = "3c3f2e5e-9ad4-413f-8618-f6ee6d42b7ff.b151e9597499e6c38ab95545bb518bc91804d098.zh-cn.xlf"
        ZhXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10659b8f251bdfd0ea5412ac80bccc419fd0af1c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3c3f2e5e-9ad4-413f-8618-f6ee6d42b7ff.b151e9597499e6c38ab95545bb518bc91804d098.zh-cn.xlf"
        ZhHandoffDt  = "2016-03-12 18:31:04"
        DeXlfFile    = "3c3f2e5e-9ad4-413f-8618-f6ee6d42b7ff.b151e9597499e6c38ab95545bb518bc91804d098.de-de.xlf"
        DeXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/331dd98359ab1a989554166435b14cd1b3b95a57/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3c3f2e5e-9ad4-413f-8618-f6ee6d42b7ff.b151e9597499e6c38ab95545bb518bc91804d098.de-de.xlf"
        DeHandoffDt  = "2016-03-12 18:31:07"
    }
)

$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$row = 4
foreach ($rec in $records) {
    $wsOverview.Cells.Item($row, 1).Value = ($rec.Uuid + ".md")
    $wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($row, 1), $rec.MdUrl, "", "", ($rec.Uuid + ".md"))
    $wsOverview.Cells.Item($row, 2).Value = $statusReady
    $wsOverview.Cells.Item($row, 3).Value = $statusReady
    $wsOverview.Cells.Item($row, 4).Value = $rec.HandoffDate
    $row++
}

# ---------------------------------------------------------------------
# Per-language detail sheets:
# Source File Name | File Extension | Status | Latest Handoff File |
# Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; XlfFileKey = "ZhXlfFile"; XlfUrlKey = "ZhXlfUrl"; HandoffDtKey = "ZhHandoffDt" },
    @{ Name = "de-de"; XlfFileKey = "DeXlfFile"; XlfUrlKey = "DeXlfUrl"; HandoffDtKey = "DeHandoffDt" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $row = 4
    foreach ($rec in $records) {
        $mdName = $rec.Uuid + ".md"
        $xlfFile = $rec[$lang.XlfFileKey]
        $xlfUrl = $rec[$lang.XlfUrlKey]
        $handoffDt = $rec[$lang.HandoffDtKey]

        # A: Source File Name (hyperlink to the .md source)
        $ws.Cells.Item($row, 1).Value = $mdName
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $rec.MdUrl, "", "", $mdName)

        # B: File Extension (hyperlink, displayed as ".md")
        $ws.Cells.Item($row, 2).Value = ".md"
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $rec.MdUrl, "", "", ".md")

        # C: Status
        $ws.Cells.Item($row, 3).Value = $statusReady

        # D: Latest Handoff File (hyperlink to the .xlf handoff file)
        $ws.Cells.Item($row, 4).Value = $xlfFile
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 4), $xlfUrl, "", "", $xlfFile)

        # E: Latest Handoff Datetime
        $ws.Cells.Item($row, 5).Value = $handoffDt
        $ws.Cells.Item($row, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

        # F, G: Latest Target File / Latest Handback File -- left blank,
        # this record has not been handed back yet.

        # H: Latest Handback DateTime (never handed back -> sentinel date)
        $ws.Cells.Item($row, 8).Value = "0001-01-01 00:00:00"

        # I: Handoff Reason
        $ws.Cells.Item($row, 9).Value = "Include"

        # J, K: Dependency From / Error Detail -- left blank.

        $row++
    }
}
